# Regenerate save_data column G ("K") with freshly calculated strikeout values.
# (regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new value for column G (7), i.e. the "K" column.
$kValues = @{
    2  = 1
    3  = 0
    4  = 0
    5  = 2
    6  = 2
    7  = 3
    8  = 1
    9  = 1
    10 = 1
    11 = 1
    12 = 1
    13 = 0
    14 = 3
    15 = 0
    16 = 0
    18 = 1
    19 = 0
    20 = 1
    21 = 1
    22 = 0
    23 = 1
    24 = 1
    25 = 1
    26 = 2
    27 = 0
    28 = 0
    29 = 1
    30 = 0
    31 = 3
    32 = 1
    33 = 1
    34 = 2
    35 = 1
    36 = 2
    38 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
